# Apply updated "想去人数" (want-to-go count) values scraped at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 1527
$ws1.Range("F13").Value = 684
$ws1.Range("F14").Value = 1847
$ws1.Range("F23").Value = 1301
$ws1.Range("F24").Value = 425
$ws1.Range("F25").Value = 519
$ws1.Range("F26").Value = 197
$ws1.Range("F27").Value = 6846
$ws1.Range("F28").Value = 5617

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 27

# --- Sheet "本地生活" --- (no changes)

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 27
$ws4.Range("F5").Value = 0
$ws4.Range("F14").Value = 1527
$ws4.Range("F16").Value = 684
$ws4.Range("F17").Value = 1847
$ws4.Range("F28").Value = 1301
$ws4.Range("F29").Value = 425
$ws4.Range("F30").Value = 519
$ws4.Range("F31").Value = 197
$ws4.Range("F32").Value = 6846
$ws4.Range("F33").Value = 5618
